# Ergebnisse vom 2. Spieltag ergänzt
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gesamttabelle")
$ws.Activate()

# Row 6: 2. Spieltag / Olympiade / Erspielte Punkte
$ws.Range("D6").Value = 26
$ws.Range("E6").Value = 14
$ws.Range("F6").Value = 18
$ws.Range("G6").Value = 11
$ws.Range("H6").Value = 15
$ws.Range("I6").Value = 0

# Row 7: 2. Spieltag / Olympiade / Würfelbonus
$ws.Range("D7").Value = 5
$ws.Range("F7").Value = 9
$ws.Range("G7").Value = 4
$ws.Range("H7").Value = 8

# Update selection to match the saved workbook state
$ws.Range("I7").Select()
